# Module 5 && Module 4 clean up
# Replace the placeholder "0" values (rows 2-26) with the real list of
# reddit post permalinks (rows 2-11), then drop the now-unused trailing
# rows so the sheet's used range shrinks from A1:A26 down to A1:A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 12-26 entirely (they held placeholder zeros and the old,
# now-replaced, link list that used to live in rows 23-26).
$ws.Rows("12:26").Delete()

# New permalinks that belong in A2:A11.
$values = @(
  "/r/EarthPorn/comments/g7nxna/northern_lights_dancing_over_the_tailrace_in/",
  "/r/EarthPorn/comments/g7soyl/spooky_fog_during_hike_to_kendall_knob_washington/",
  "/r/goodnews/comments/fk33b0/modmessage_visit_the_good_news_discord_server/",
  "/r/goodnews/comments/g783r7/good_news_its_friday_whats_your_feelgood_story/",
  "/r/learnpython/comments/g4iiwc/ask_anything_monday_weekly_thread/",
  "/r/learnpython/comments/g7rpwu/ok_so_im_committed_to_1_year_of_coding_in_python/",
  "/r/pics/comments/fjn0j9/important_psa_no_you_did_not_win_a_gift_card/",
  "/r/pics/comments/g7rlvg/when_a_cat_runs_to_the_fridge_every_time_it_opens/",
  "/r/Python/comments/g5fwr9/whats_everyone_working_on_this_week/",
  "/r/Python/comments/g7q2ej/my_pothole_detector_used_yolov3_annotated_images/"
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Match the saved view state: selection parked on I50.
$ws.Range("I50").Select()

Write-Output "Updated reddit_posts Sheet1 with $($values.Length) links (A2:A11)."
